# Apply "cement and other calibration adjustments" edits to the "energy" sheet:
# add a parameter/source reference table in columns J:K, rows 34-50.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("energy")

# Widen column J to fit the new "parameter" labels.
$ws.Columns.Item(10).ColumnWidth = 23.75

# Populate J34:K50 in the same order the values were originally entered
# (first-use order determines shared-string ids, so keep this sequence).
$ws.Range("J35").Value = "historical floorspace"
$ws.Range("J36").Value = "satiation of floorspace"
$ws.Range("J39").Value = "satiation of services"
$ws.Range("K34").Value = "source"
$ws.Range("K35").Value = "processed from other sources"
$ws.Range("K36").Value = "assumed"
$ws.Range("K37").Value = "assumed"
$ws.Range("K44").Value = "assumed"
$ws.Range("J37").Value = "price elasticity of floorspace"
$ws.Range("J38").Value = "degree days"
$ws.Range("K38").Value = "processed by Yuyu (GIS)"
$ws.Range("K39").Value = "processed from base year data, multiplied by exog multipliers and other adjustments"
$ws.Range("K40").Value = "non-US regions: heating and cooling satiation levels equal to US base-year level, multiplied by HDD or CDD ratio"
$ws.Range("J41").Value = "internal gain scaler"
$ws.Range("K41").Value = "exogenous (set to artificially low level in order to nullify potentially negative impacts of this parameter)"
$ws.Range("J34").Value = "parameter"
$ws.Range("J42").Value = "service output (`"base service`")"
$ws.Range("K42").Value = "aggregated as the sum of all service outputs (at the tech level)"
$ws.Range("J43").Value = "shell conductance"
$ws.Range("K43").Value = "USA: US input set. Other regions: assigned exogenous tech change"
$ws.Range("J44").Value = "fuel preference elasticity"
$ws.Range("J45").Value = "efficiency"
$ws.Range("K45").Value = "usa from detailed US data. Other regions multiplied by exogenous region- and time-specific adjustment factors"
$ws.Range("J46").Value = "callibrated energy consumption"
$ws.Range("K46").Value = "IEA energy balances' building energy demand by fuel, multiplied by service allocation shares"
$ws.Range("J47").Value = "service allocation shares"
$ws.Range("K47").Value = "requires processing; from specific data sources for countries and regions"
$ws.Range("J48").Value = "non-fuel costs"
$ws.Range("K48").Value = "calculated from exogenous assumptions about capital costs, interest rates, O&M costs, UECs, and efficiencies of each technology"
$ws.Range("J49").Value = "internal load fraction"
$ws.Range("K49").Value = "exogenous (equal for all techs right now)"
$ws.Range("J50").Value = "retirement"
$ws.Range("K50").Value = "exogenous (from USA data processing)"

# Restore the view: scrolled down so row 26 is at the top, with J41 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J41").Select()
